# Mifos -> Finflux: insert a new (blank) column into the "Repayment schedule"
# sheet immediately before what used to be column N ("Late").
#
# Before: ... M=In Advance | N=Late   | O=(blank) | P=Outstanding
# After:  ... M=In Advance | N=(blank)| O=Late    | P=(blank) | Q=Outstanding
#
# i.e. a whole new blank column is inserted at N, shifting the old N/O/P
# columns one place to the right (N->O, O->P, P->Q).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Inserting a whole column at N shifts existing N:P data right by one,
# exactly matching the diff (old N becomes O, old O becomes P, old P
# becomes Q), and leaves the freshly inserted column N blank.
[void]$ws.Columns("N:N").Insert()

# The author's cursor ended up on K18 afterwards (was K9 before the edit).
[void]$ws.Range("K18").Select()
